# Auto-generated Excel COM-interop script
# Updates cached numeric values (H:N columns) across 8 job sheets,
# mirroring a scheduled-runner refresh of market-price driven profit figures.

$wb = $excel.ActiveWorkbook

# ALC row 132: Fast-forwarding Flora
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 4190.875
$ws.Range("I132").Value = 4240.0356
$ws.Range("J132").Value = 3846.75
$ws.Range("K132").Value = 12720.1068
$ws.Range("L132").Value = 11540.25
$ws.Range("M132").Value = -10190.1068
$ws.Range("N132").Value = -16600.25

# ALC row 135: For Tired Minds
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 1597.7142
$ws.Range("I135").Value = 1705.3334
$ws.Range("J135").Value = 1328.6666
$ws.Range("K135").Value = 15348.0006
$ws.Range("L135").Value = 11957.9994
$ws.Range("M135").Value = -12813.0006
$ws.Range("N135").Value = -17027.9994

# ALC row 137: Cutting Edge of Culinary Quality
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1327.72
$ws.Range("I137").Value = 1193.125
$ws.Range("K137").Value = 3579.375
$ws.Range("M137").Value = -1029.375

# ALC row 138: All-night Crafting
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2516.7869
$ws.Range("J138").Value = 2420
$ws.Range("L138").Value = 7260
$ws.Range("N138").Value = -17540

# ALC row 141: Remedy for Reason
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 5001.92
$ws.Range("I141").Value = 2092
$ws.Range("J141").Value = 11185.5
$ws.Range("K141").Value = 6276
$ws.Range("L141").Value = 33556.5
$ws.Range("M141").Value = -1096
$ws.Range("N141").Value = -43916.5

# ARM row 32: Ingot We Trust
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 383078.06
$ws.Range("I32").Value = 440723.97
$ws.Range("K32").Value = 440723.97
$ws.Range("M32").Value = -440436.97

# ARM row 61: Dealing with the Tough Stuff
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 13336022
$ws.Range("I61").Value = 55556990
$ws.Range("J61").Value = 3083.7896
$ws.Range("K61").Value = 55556990
$ws.Range("L61").Value = 3083.7896
$ws.Range("M61").Value = -55556778
$ws.Range("N61").Value = -3507.7896

# ARM row 74: As the Bolt Flies
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1469.2778
$ws.Range("I74").Value = 1038.3334
$ws.Range("J74").Value = 1684.75
$ws.Range("K74").Value = 1038.3334
$ws.Range("L74").Value = 1684.75
$ws.Range("M74").Value = -164.3334
$ws.Range("N74").Value = -3432.75

# ARM row 77: Heavy Metal Banned (L)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1469.2778
$ws.Range("I77").Value = 1038.3334
$ws.Range("J77").Value = 1684.75
$ws.Range("K77").Value = 5191.666999999999
$ws.Range("L77").Value = 8423.75
$ws.Range("M77").Value = -823.6669999999995
$ws.Range("N77").Value = -17159.75

# ARM row 122: Haste for High Durium
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1424.5555
$ws.Range("I122").Value = 1137
$ws.Range("J122").Value = 1999.6666
$ws.Range("K122").Value = 3411
$ws.Range("L122").Value = 5998.9998
$ws.Range("M122").Value = -961
$ws.Range("N122").Value = -10898.9998

# ARM row 136: Metal with Mettle
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 13336022
$ws.Range("I136").Value = 55556990
$ws.Range("J136").Value = 3083.7896
$ws.Range("K136").Value = 166670970
$ws.Range("L136").Value = 9251.3688
$ws.Range("M136").Value = -166668420
$ws.Range("N136").Value = -14351.3688

# CRP row 31: Wall Not Found
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3561.7964
$ws.Range("I31").Value = 1352.5264
$ws.Range("J31").Value = 4761.1143
$ws.Range("K31").Value = 1352.5264
$ws.Range("L31").Value = 4761.1143
$ws.Range("M31").Value = -1057.5264
$ws.Range("N31").Value = -5351.1143

# CRP row 34: Armoires of the Rich and Famous
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3561.7964
$ws.Range("I34").Value = 1352.5264
$ws.Range("J34").Value = 4761.1143
$ws.Range("K34").Value = 1352.5264
$ws.Range("L34").Value = 4761.1143
$ws.Range("M34").Value = -1150.5264
$ws.Range("N34").Value = -5165.1143

# CRP row 58: You Do the Heavy Lifting
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2905.6428
$ws.Range("I58").Value = 3022.375
$ws.Range("J58").Value = 2750
$ws.Range("K58").Value = 3022.375
$ws.Range("L58").Value = 2750
$ws.Range("M58").Value = -2819.375
$ws.Range("N58").Value = -3156

# CRP row 122: Timber of Tenkonto
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1789.0769
$ws.Range("I122").Value = 1514.7142
$ws.Range("J122").Value = 1890.1578
$ws.Range("K122").Value = 4544.142599999999
$ws.Range("L122").Value = 5670.4734
$ws.Range("M122").Value = -2094.142599999999
$ws.Range("N122").Value = -10570.4734

# CRP row 132: Hull Lotta Damage
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 15154329
$ws.Range("I132").Value = 1499.6666
$ws.Range("J132").Value = 20836640
$ws.Range("K132").Value = 4498.9998
$ws.Range("L132").Value = 62509920
$ws.Range("M132").Value = -1968.9998
$ws.Range("N132").Value = -62514980

# CRP row 134: Wood You Be Quiet
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 994.03845
$ws.Range("I134").Value = 737.7727
$ws.Range("K134").Value = 2213.3181
$ws.Range("M134").Value = 321.6819

# CRP row 136: Turali Quality
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2905.6428
$ws.Range("I136").Value = 3022.375
$ws.Range("J136").Value = 2750
$ws.Range("K136").Value = 9067.125
$ws.Range("L136").Value = 8250
$ws.Range("M136").Value = -6517.125
$ws.Range("N136").Value = -13350

# CUL row 12: Butter Me Up
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 102.20588
$ws.Range("I12").Value = 2.875
$ws.Range("J12").Value = 132.76923
$ws.Range("K12").Value = 8.625
$ws.Range("L12").Value = 398.30769
$ws.Range("M12").Value = 164.375
$ws.Range("N12").Value = -744.30769

# CUL row 68: Such a Butter Face
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2290.0618
$ws.Range("I68").Value = 3143.5945
$ws.Range("J68").Value = 1763.7167
$ws.Range("K68").Value = 9430.783500000001
$ws.Range("L68").Value = 5291.1501
$ws.Range("M68").Value = -8619.783500000001
$ws.Range("N68").Value = -6913.1501

# CUL row 71: No Margarine of Error (L)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 2290.0618
$ws.Range("I71").Value = 3143.5945
$ws.Range("J71").Value = 1763.7167
$ws.Range("K71").Value = 28292.3505
$ws.Range("L71").Value = 15873.4503
$ws.Range("M71").Value = -24236.3505
$ws.Range("N71").Value = -23985.4503

# CUL row 113: Can't Eat Just One
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1074.9667
$ws.Range("I113").Value = 581.93335
$ws.Range("J113").Value = 1568
$ws.Range("K113").Value = 1745.80005
$ws.Range("L113").Value = 4704
$ws.Range("M113").Value = 424.1999499999999
$ws.Range("N113").Value = -9044

# CUL row 121: A Cookie for Your Troubles
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 963.88525
$ws.Range("I121").Value = 511.9
$ws.Range("J121").Value = 1052.5098
$ws.Range("K121").Value = 1535.7
$ws.Range("L121").Value = 3157.5294
$ws.Range("M121").Value = -225.6999999999998
$ws.Range("N121").Value = -5777.5294

# CUL row 131: The Mountain Steeped
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 953.8095
$ws.Range("J131").Value = 1184
$ws.Range("L131").Value = 3552
$ws.Range("N131").Value = -13632

# GSM row 82: Appeasing the Astromancer
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H82").Value = 25999.646
$ws.Range("J82").Value = 25999.646
$ws.Range("L82").Value = 25999.646
$ws.Range("N82").Value = -26765.646

# GSM row 85: Silver Bar of Upcycling (L)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H85").Value = 25999.646
$ws.Range("J85").Value = 25999.646
$ws.Range("L85").Value = 25999.646
$ws.Range("N85").Value = -28651.646

# GSM row 122: Awarding Academic Excellence
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4716.3335
$ws.Range("I122").Value = 3989.25
$ws.Range("K122").Value = 11967.75
$ws.Range("M122").Value = -9517.75

# LTW row 136: Respect for Br'aax
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2688953
$ws.Range("I136").Value = 933.63336
$ws.Range("J136").Value = 5208971.5
$ws.Range("K136").Value = 2800.90008
$ws.Range("L136").Value = 15626914.5
$ws.Range("M136").Value = -250.9000800000003
$ws.Range("N136").Value = -15632014.5

# WVR row 132: Comfy Cabins
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 5210204.5
$ws.Range("I132").Value = 1752.875
$ws.Range("J132").Value = 12154807
$ws.Range("K132").Value = 5258.625
$ws.Range("L132").Value = 36464421
$ws.Range("M132").Value = -2728.625
$ws.Range("N132").Value = -36469481

# WVR row 136: Weaving the Envelope
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2287.4333
$ws.Range("I136").Value = 2094.3777
$ws.Range("J136").Value = 2866.6
$ws.Range("K136").Value = 6283.1331
$ws.Range("L136").Value = 8599.799999999999
$ws.Range("M136").Value = -3733.1331
$ws.Range("N136").Value = -13699.8
